# A new weekly price record (row 161, dated 45075 = 2023-05-29) is inserted
# into the "Frambuesa" price list. All following rows (old 161-167) shift
# down by one (new 162-168); their contents are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 161 (and everything below it) down by one row,
# opening up a blank row 161 for the new record.
$ws.Rows(161).Insert()

# Populate the newly inserted row 161 with the new record's data.
$ws.Range("A161").Value = 9
$ws.Range("B161").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C161").Value = "Metropolitana"
$ws.Range("D161").Value = 45075
$ws.Range("E161").Value = 13
$ws.Range("F161").Value = "Fruta"
$ws.Range("G161").Value = 100101
$ws.Range("H161").Value = "Berries"
$ws.Range("I161").Value = 100101004
$ws.Range("J161").Value = "Frambuesa"
$ws.Range("K161").Value = "Sin especificar"
$ws.Range("L161").Value = "Primera"
$ws.Range("M161").Value = 350
$ws.Range("N161").Value = 9000
$ws.Range("O161").Value = 9500
$ws.Range("P161").Value = 9286
$ws.Range("Q161").Value = "$/bandeja 2 kilos"
$ws.Range("R161").Value = "Provincia de Linares"
$ws.Range("S161").Value = 4643
$ws.Range("T161").Value = 2
